$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 28 (SC 92) first, then row 26 (RM 232), bottom-up so row indices
# for the earlier deletion remain valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# After the two deletions, the remaining rows have shifted up:
#   original row 27 (SC 5)    -> now row 26
#   original row 29 (SC 101)  -> now row 27
#   original row 35 (SC 232)  -> now row 33
# Update column F (imputed values) for these rows to match the new data.
$ws.Range("F26").Value = 17.38
$ws.Range("F27").ClearContents()
$ws.Range("F33").Value = 17.53
